# Insert a new empty paragraph right after the "#endif" paragraph.
$d = $word.ActiveDocument

# Locate the paragraph whose text is "#endif" (the end of the include guard).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "#endif") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
}
